$wb = $excel.ActiveWorkbook

# --- Sheet: "Data Set 0 Timings (Pd)" (rows 2-6) ---
$wsPd = $wb.Worksheets.Item("Data Set 0 Timings (Pd)")
$wsPd.Range("F2").Value = 2.088
$wsPd.Range("H2").Value = 1.27
$wsPd.Range("I2").Value = 0.376
$wsPd.Range("J2").Value = 0.369
$wsPd.Range("L2").Value = 0.0002

$wsPd.Range("F3").Value = 2.053
$wsPd.Range("G3").Value = 0.051
$wsPd.Range("H3").Value = 1.117
$wsPd.Range("I3").Value = 0.369
$wsPd.Range("J3").Value = 0.434
$wsPd.Range("L3").Value = 0.00018

$wsPd.Range("F4").Value = 2.023
$wsPd.Range("G4").Value = 0.041
$wsPd.Range("H4").Value = 1.131
$wsPd.Range("I4").Value = 0.44
$wsPd.Range("J4").Value = 0.381
$wsPd.Range("L4").Value = 0.00019

$wsPd.Range("F5").Value = 2.058
$wsPd.Range("G5").Value = 0.043
$wsPd.Range("H5").Value = 1.217
$wsPd.Range("I5").Value = 0.378
$wsPd.Range("J5").Value = 0.39
$wsPd.Range("L5").Value = 0.0002

$wsPd.Range("F6").Value = 2.178
$wsPd.Range("G6").Value = 0.043
$wsPd.Range("H6").Value = 1.253
$wsPd.Range("I6").Value = 0.398
$wsPd.Range("J6").Value = 0.455
$wsPd.Range("L6").Value = 0.0002

# --- Sheet: "Data Set 0 Timings (TD)" (rows 2-6) ---
$wsTd = $wb.Worksheets.Item("Data Set 0 Timings (TD)")
$wsTd.Range("F2").Value = 3.114
$wsTd.Range("I2").Value = 1.097
$wsTd.Range("J2").Value = 1.225
$wsTd.Range("L2").Value = 0.00022

$wsTd.Range("F3").Value = 3.191
$wsTd.Range("H3").Value = 0.004
$wsTd.Range("I3").Value = 1.235
$wsTd.Range("J3").Value = 1.216
$wsTd.Range("L3").Value = 0.00023

$wsTd.Range("F4").Value = 3.118
$wsTd.Range("I4").Value = 1.181
$wsTd.Range("J4").Value = 1.207
$wsTd.Range("L4").Value = 0.00022

$wsTd.Range("F5").Value = 3.165
$wsTd.Range("H5").Value = 0.004
$wsTd.Range("I5").Value = 1.212
$wsTd.Range("J5").Value = 1.219
$wsTd.Range("L5").Value = 0.00022

$wsTd.Range("F6").Value = 3.233
$wsTd.Range("H6").Value = 0.004
$wsTd.Range("I6").Value = 1.254
$wsTd.Range("J6").Value = 1.247

# --- Sheet: "Data Set 0 Timings (combined)" (rows 2-6 mirror Pd, rows 7-11 mirror TD) ---
$wsComb = $wb.Worksheets.Item("Data Set 0 Timings (combined)")
$wsComb.Range("F2").Value = 2.088
$wsComb.Range("H2").Value = 1.27
$wsComb.Range("I2").Value = 0.376
$wsComb.Range("J2").Value = 0.369
$wsComb.Range("L2").Value = 0.0002

$wsComb.Range("F3").Value = 2.053
$wsComb.Range("G3").Value = 0.051
$wsComb.Range("H3").Value = 1.117
$wsComb.Range("I3").Value = 0.369
$wsComb.Range("J3").Value = 0.434
$wsComb.Range("L3").Value = 0.00018

$wsComb.Range("F4").Value = 2.023
$wsComb.Range("G4").Value = 0.041
$wsComb.Range("H4").Value = 1.131
$wsComb.Range("I4").Value = 0.44
$wsComb.Range("J4").Value = 0.381
$wsComb.Range("L4").Value = 0.00019

$wsComb.Range("F5").Value = 2.058
$wsComb.Range("G5").Value = 0.043
$wsComb.Range("H5").Value = 1.217
$wsComb.Range("I5").Value = 0.378
$wsComb.Range("J5").Value = 0.39
$wsComb.Range("L5").Value = 0.0002

$wsComb.Range("F6").Value = 2.178
$wsComb.Range("G6").Value = 0.043
$wsComb.Range("H6").Value = 1.253
$wsComb.Range("I6").Value = 0.398
$wsComb.Range("J6").Value = 0.455
$wsComb.Range("L6").Value = 0.0002

$wsComb.Range("F7").Value = 3.114
$wsComb.Range("I7").Value = 1.097
$wsComb.Range("J7").Value = 1.225
$wsComb.Range("L7").Value = 0.00022

$wsComb.Range("F8").Value = 3.191
$wsComb.Range("H8").Value = 0.004
$wsComb.Range("I8").Value = 1.235
$wsComb.Range("J8").Value = 1.216
$wsComb.Range("L8").Value = 0.00023

$wsComb.Range("F9").Value = 3.118
$wsComb.Range("I9").Value = 1.181
$wsComb.Range("J9").Value = 1.207
$wsComb.Range("L9").Value = 0.00022

$wsComb.Range("F10").Value = 3.165
$wsComb.Range("H10").Value = 0.004
$wsComb.Range("I10").Value = 1.212
$wsComb.Range("J10").Value = 1.219
$wsComb.Range("L10").Value = 0.00022

$wsComb.Range("F11").Value = 3.233
$wsComb.Range("H11").Value = 0.004
$wsComb.Range("I11").Value = 1.254
$wsComb.Range("J11").Value = 1.247
